# Weekly refresh of the Camote (Vega Modelo de Temuco) price sheet:
# a new daily record is inserted as row 94, pushing the previously
# existing rows 94-172 down to 95-173 (dimension grows from R172 to R173).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94 - shifts rows 94:172 down to 95:173
# and keeps the D-column date formatting.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new record.
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 44978
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 100114002
$ws.Range("G94").Value = "Camote"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 10
$ws.Range("K94").Value = 26000
$ws.Range("L94").Value = 26000
$ws.Range("M94").Value = 26000
$ws.Range("N94").Value = "$/malla 20 kilos"
$ws.Range("O94").Value = "Perú"
$ws.Range("P94").Value = 1300
$ws.Range("Q94").Value = 20
$ws.Range("R94").Value = "Hortaliza"
